$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fn1"
$ws.Cells.Item(2,3).Value = "Mag"
$ws.Cells.Item(2,4).Value = "M2"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 19.95578266666667
$ws.Cells.Item(2,8).Value = 59.867348
$ws.Cells.Item(2,9).Value = 0.0117373419656925
$ws.Cells.Item(2,10).Value = 0.0117373419656925
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.7885686666666668
$ws.Cells.Item(2,14).Value = 2.365706
$ws.Cells.Item(2,15).Value = 0.4566863346753138
$ws.Cells.Item(2,16).Value = 0.4566863346753137
$ws.Cells.Item(2,17).Value = 15.73650492974311
$ws.Cells.Item(2,18).Value = 141.628544367688
$ws.Cells.Item(2,19).Value = 0.005360283681142849
$ws.Cells.Item(2,20).Value = 0.005360283681142848

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fn1"
$ws.Cells.Item(3,3).Value = "Mag"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 19.95578266666667
$ws.Cells.Item(3,8).Value = 59.867348
$ws.Cells.Item(3,9).Value = 0.0117373419656925
$ws.Cells.Item(3,10).Value = 0.0117373419656925
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.9381496666666668
$ws.Cells.Item(3,14).Value = 2.814449
$ws.Cells.Item(3,15).Value = 0.5433136653246862
$ws.Cells.Item(3,16).Value = 0.5433136653246862
$ws.Cells.Item(3,17).Value = 18.72151085680578
$ws.Cells.Item(3,18).Value = 168.493597711252
$ws.Cells.Item(3,19).Value = 0.006377058284549647
$ws.Cells.Item(3,20).Value = 0.006377058284549647

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Fn1"
$ws.Cells.Item(4,3).Value = "Mag"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1637.343343333333
$ws.Cells.Item(4,8).Value = 4912.03003
$ws.Cells.Item(4,9).Value = 0.9630320723052701
$ws.Cells.Item(4,10).Value = 0.9630320723052702
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.7885686666666668
$ws.Cells.Item(4,14).Value = 2.365706
$ws.Cells.Item(4,15).Value = 0.4566863346753138
$ws.Cells.Item(4,16).Value = 0.4566863346753137
$ws.Cells.Item(4,17).Value = 1291.157657127909
$ws.Cells.Item(4,18).Value = 11620.41891415118
$ws.Cells.Item(4,19).Value = 0.4398035872758655
$ws.Cells.Item(4,20).Value = 0.4398035872758655

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fn1"
$ws.Cells.Item(5,3).Value = "Mag"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1637.343343333333
$ws.Cells.Item(5,8).Value = 4912.03003
$ws.Cells.Item(5,9).Value = 0.9630320723052701
$ws.Cells.Item(5,10).Value = 0.9630320723052702
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.9381496666666668
$ws.Cells.Item(5,14).Value = 2.814449
$ws.Cells.Item(5,15).Value = 0.5433136653246862
$ws.Cells.Item(5,16).Value = 0.5433136653246862
$ws.Cells.Item(5,17).Value = 1536.073111767052
$ws.Cells.Item(5,18).Value = 13824.65800590347
$ws.Cells.Item(5,19).Value = 0.5232284850294044
$ws.Cells.Item(5,20).Value = 0.5232284850294046

# Row 6
$ws.Cells.Item(6,1).Value = "M2"
$ws.Cells.Item(6,2).Value = "Fn1"
$ws.Cells.Item(6,3).Value = "Mag"
$ws.Cells.Item(6,4).Value = "M2"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 17.50081933333334
$ws.Cells.Item(6,8).Value = 52.502458
$ws.Cells.Item(6,9).Value = 0.01029341242216722
$ws.Cells.Item(6,10).Value = 0.01029341242216722
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.7885686666666668
$ws.Cells.Item(6,14).Value = 2.365706
$ws.Cells.Item(6,15).Value = 0.4566863346753138
$ws.Cells.Item(6,16).Value = 0.4566863346753137
$ws.Cells.Item(6,17).Value = 13.80059776726089
$ws.Cells.Item(6,18).Value = 124.205379905348
$ws.Cells.Item(6,19).Value = 0.00470086079038089
$ws.Cells.Item(6,20).Value = 0.00470086079038089

# Row 7
$ws.Cells.Item(7,1).Value = "M2"
$ws.Cells.Item(7,2).Value = "Fn1"
$ws.Cells.Item(7,3).Value = "Mag"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 17.50081933333334
$ws.Cells.Item(7,8).Value = 52.502458
$ws.Cells.Item(7,9).Value = 0.01029341242216722
$ws.Cells.Item(7,10).Value = 0.01029341242216722
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.9381496666666668
$ws.Cells.Item(7,14).Value = 2.814449
$ws.Cells.Item(7,15).Value = 0.5433136653246862
$ws.Cells.Item(7,16).Value = 0.5433136653246862
$ws.Cells.Item(7,17).Value = 16.41838782396023
$ws.Cells.Item(7,18).Value = 147.765490415642
$ws.Cells.Item(7,19).Value = 0.005592551631786327
$ws.Cells.Item(7,20).Value = 0.005592551631786328

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Fn1"
$ws.Cells.Item(8,3).Value = "Mag"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 25.39612333333333
$ws.Cells.Item(8,8).Value = 76.18836999999999
$ws.Cells.Item(8,9).Value = 0.01493717330687017
$ws.Cells.Item(8,10).Value = 0.01493717330687017
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.7885686666666668
$ws.Cells.Item(8,14).Value = 2.365706
$ws.Cells.Item(8,15).Value = 0.4566863346753138
$ws.Cells.Item(8,16).Value = 0.4566863346753137
$ws.Cells.Item(8,17).Value = 20.02658711546889
$ws.Cells.Item(8,18).Value = 180.23928403922
$ws.Cells.Item(8,19).Value = 0.006821602927924473
$ws.Cells.Item(8,20).Value = 0.006821602927924472

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Fn1"
$ws.Cells.Item(9,3).Value = "Mag"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 25.39612333333333
$ws.Cells.Item(9,8).Value = 76.18836999999999
$ws.Cells.Item(9,9).Value = 0.01493717330687017
$ws.Cells.Item(9,10).Value = 0.01493717330687017
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.9381496666666668
$ws.Cells.Item(9,14).Value = 2.814449
$ws.Cells.Item(9,15).Value = 0.5433136653246862
$ws.Cells.Item(9,16).Value = 0.5433136653246862
$ws.Cells.Item(9,17).Value = 23.82536463979222
$ws.Cells.Item(9,18).Value = 214.42828175813
$ws.Cells.Item(9,19).Value = 0.008115570378945695
$ws.Cells.Item(9,20).Value = 0.008115570378945695
